$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Test 039): mark implemented, update notes
$ws.Range("D7").Value = "Y"
$ws.Range("E7").Value = "Reworked with a tree for variable usages."

# Row 11 (Test 048): mark implemented
$ws.Range("D11").Value = "Y"

# Row 12 (Test 049): mark implemented
$ws.Range("D12").Value = "Y"

# Row 18 (Test 156): mark implemented, add note
$ws.Range("D18").Value = "Y"
$ws.Range("E18").Value = "Tests pass"

# Row 19 (Test 157): mark implemented, add note
$ws.Range("D19").Value = "Y"
$ws.Range("E19").Value = "Tests pass"

# Row 13 (Test 053): mark implemented, add note
$ws.Range("D13").Value = "Y"
$ws.Range("E13").Value = "Sets value to INT_MAX"

# Row 14 (Test 054): mark implemented, add note
$ws.Range("D14").Value = "Y"
$ws.Range("E14").Value = "Sets value to INT_MIN"

# C16's style picks up a text number format (quote-prefixed cell)
$ws.Range("C16").NumberFormat = "@"

# Update selection to reflect last active cell
$ws.Range("C16").Select()
